$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 1")

# Add a trailing period to each requirement description in column B (rows 2-17),
# if it isn't already terminated with one.
for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $text = $cell.Value2
    if ($text -ne $null -and -not $text.EndsWith(".")) {
        $cell.Value2 = $text + "."
    }
}

# Leave the selection on B13, matching the last-edited cell in the sheet view.
$ws.Range("B13").Select()
